# Add connection to database: new "Sheet2" with the same comparison table
# but refreshed metrics (computed "using database" instead of "using csv file"),
# and mark the original sheet with a "using csv file" note.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: note how the data was produced ---
$ws1.Range("A9").Value = "using csv file"
$ws1.PageSetup.Orientation = 1
$ws1.Range("A1:E7").Select() | Out-Null

# --- New Sheet2: same layout, refreshed numbers from the database run ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Comparison of Models"

$ws2.Range("B2").Value = "Accuracy"
$ws2.Range("C2").Value = "Precision:"
$ws2.Range("D2").Value = "Recall:"
$ws2.Range("E2").Value = "F1"

$ws2.Range("A3").Value = "Random Forest"
$ws2.Range("B3").Value = 0.803
$ws2.Range("C3").Value = 0.39
$ws2.Range("D3").Value = 0.24
$ws2.Range("E3").Value = 0.29

$ws2.Range("A4").Value = "Logistic Regression"
$ws2.Range("B4").Value = 0.83
$ws2.Range("C4").Value = 0.57
$ws2.Range("D4").Value = 0.09
$ws2.Range("E4").Value = 0.15

$ws2.Range("A5").Value = "Support Vector Machine"
$ws2.Range("B5").Value = 0.831
$ws2.Range("C5").Value = 0.61
$ws2.Range("D5").Value = 0.07
$ws2.Range("E5").Value = 0.12

$ws2.Range("A6").Value = "Deep Learning"
$ws2.Range("B6").Value = 0.831
$ws2.Range("C6").Value = 0.58
$ws2.Range("D6").Value = 0.09
$ws2.Range("E6").Value = 0.16

$ws2.Range("A7").Value = "Deep Learning Final"
$ws2.Range("B7").Value = 0.832
$ws2.Range("C7").Value = 0.58
$ws2.Range("D7").Value = 0.1
$ws2.Range("E7").Value = 0.17

$ws2.Range("A9").Value = "using database"

$ws2.Range("B3:E7").NumberFormat = "0.000"
$ws2.Columns("A").ColumnWidth = 21.6

$ws2.Range("C7").Select() | Out-Null
